$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the active selection (A16 -> A11)
$ws.Range("A11").Select() | Out-Null

# Rewrite the data rows (A: item name, B: quantity) to the new snapshot values
$ws.Range("A2").Value = 'ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы'
$ws.Range("B2").Value = 81024
$ws.Range("A3").Value = 'Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г'
$ws.Range("B3").Value = 1595
$ws.Range("A4").Value = 'Шалфей листья 50г'
$ws.Range("B4").Value = 22792
$ws.Range("A5").Value = 'Бессмертник песчаный цветки 30г'
$ws.Range("B5").Value = 18259
$ws.Range("A6").Value = 'Дуба кора 75г'
$ws.Range("B6").Value = 45525
$ws.Range("A7").Value = 'Лен семена 100г'
$ws.Range("B7").Value = 39734
$ws.Range("A8").Value = 'Укроп пахучий плоды 50г'
$ws.Range("B8").Value = 44066
$ws.Range("A9").Value = 'Аир корневища 75г'
$ws.Range("B9").Value = 4875
$ws.Range("A10").Value = 'Валериана корневища с корнями 50г'
$ws.Range("B10").Value = 13746
$ws.Range("A11").Value = 'Ламинарии слоевища (морская капуста) 100г'
$ws.Range("B11").Value = 10314
$ws.Range("A12").Value = 'Ромашка цветки вн 50г'
$ws.Range("B12").Value = 72572
$ws.Range("A13").Value = 'Крушина кора 50г'
$ws.Range("B13").Value = 7650
$ws.Range("A14").Value = 'Мать-и-мачеха листья 35г'
$ws.Range("B14").Value = 19842
$ws.Range("A15").Value = 'Рябина плоды 50г'
$ws.Range("B15").Value = 1358
$ws.Range("A16").Value = 'Девясил корневища и корни 50г'
$ws.Range("B16").Value = 13645
$ws.Range("A17").Value = 'Сб. Грудной №4 50г'
$ws.Range("B17").Value = 28406
$ws.Range("A18").Value = 'Полынь горькая трава 50г'
$ws.Range("B18").Value = 32002
$ws.Range("A19").Value = 'Пижма цветки 75г'
$ws.Range("B19").Value = 12910
$ws.Range("A20").Value = 'Брусника листья 50г'
$ws.Range("B20").Value = 14039
$ws.Range("A21").Value = 'Зверобой трава 50г'
$ws.Range("B21").Value = 30282
$ws.Range("A22").Value = 'Череда трава 50г'
$ws.Range("B22").Value = 10962
$ws.Range("A23").Value = 'Тысячелистник трава 50г'
$ws.Range("B23").Value = 14057
$ws.Range("A24").Value = 'Шиповник плоды низковитаминные 50г'
$ws.Range("B24").Value = 33950
$ws.Range("A25").Value = 'Можжевельник плоды 50г'
$ws.Range("B25").Value = 12220
$ws.Range("A26").Value = 'Кукуруза столбики с рыльцами 40г'
$ws.Range("B26").Value = 25627
$ws.Range("A27").Value = 'Береза почки 50г'
$ws.Range("B27").Value = 18743
$ws.Range("A28").Value = 'Спорыш трава 50г'
$ws.Range("B28").Value = 14728
$ws.Range("A29").Value = 'Липа цветки 35г'
$ws.Range("B29").Value = 23666
$ws.Range("A30").Value = 'Сенна листья 50г'
$ws.Range("B30").Value = 23211
$ws.Range("A31").Value = 'Чага (березовый гриб) 50г'
$ws.Range("B31").Value = 30240
$ws.Range("A32").Value = 'Эвкалипт прутовидный листья 75г'
$ws.Range("B32").Value = 28855
$ws.Range("A33").Value = 'Сб. Фитонефрол (Урологический сбор) 50г'
$ws.Range("B33").Value = 20466
$ws.Range("A34").Value = 'Боярышник плоды 75г'
$ws.Range("B34").Value = 24592
$ws.Range("A35").Value = 'Толокнянка листья 50г'
$ws.Range("B35").Value = 8696
$ws.Range("A36").Value = 'Подорожник большой листья 50г'
$ws.Range("B36").Value = 11158
$ws.Range("A37").Value = 'Солодка корни 50г'
$ws.Range("B37").Value = 41497
$ws.Range("A38").Value = 'Алтей корни 75г'
$ws.Range("B38").Value = 7404
$ws.Range("A39").Value = 'Чистотел трава 50г'
$ws.Range("B39").Value = 22736
$ws.Range("A40").Value = 'Эрва шерстистая трава 30г'
$ws.Range("B40").Value = 18533
$ws.Range("A41").Value = 'Крапива листья 50г'
$ws.Range("B41").Value = 18177
$ws.Range("A42").Value = 'Мята перечная листья 50г'
$ws.Range("B42").Value = 31804
$ws.Range("A43").Value = 'Ноготки цветки 50г'
$ws.Range("B43").Value = 34042
$ws.Range("A44").Value = 'Багульник болотный побеги 50г'
$ws.Range("B44").Value = 19605
$ws.Range("A45").Value = 'Сб. Фитопектол №1 (Грудной сбор №1) 35г'
$ws.Range("B45").Value = 7069
$ws.Range("A46").Value = 'Сб. Фитопектол №2 (Грудной сбор №2) 35г'
$ws.Range("B46").Value = 9514
$ws.Range("A47").Value = 'Чабрец трава 50г'
$ws.Range("B47").Value = 36567
$ws.Range("A48").Value = 'Пустырник трава 50г'
$ws.Range("B48").Value = 28796
$ws.Range("A49").Value = 'Фп Фиточай "Лактафитол" (БАД) 20х1,5 г'
$ws.Range("B49").NumberFormat = "#,##0"
$ws.Range("B49").Value = 13935
$ws.Range("A50").Value = 'Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г'
$ws.Range("B50").Value = 3440
$ws.Range("A51").Value = 'Фп Детский травяной чай "ФармаЦветик® для иммунитета" 20х1,5 г'
$ws.Range("B51").Value = 3006
$ws.Range("A52").Value = 'Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г'
$ws.Range("B52").Value = 5450
$ws.Range("A53").Value = 'Фп Детский травяной чай "ФармаЦветик® для спокойного сна" 20х1,5 г'
$ws.Range("B53").Value = 8108
$ws.Range("A54").Value = 'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем"(БАД) 20*1,5г'
$ws.Range("B54").Value = 7170
$ws.Range("A55").Value = 'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем" (БАД) 20*1,5г'
$ws.Range("B55").Value = 8870
$ws.Range("A56").Value = 'Фп "Щедрость природы® Фиточай для иммунитета" 20х2,0 г'
$ws.Range("B56").Value = 378
$ws.Range("A57").Value = 'Фп "Щедрость природы® Фиточай диабетический" 20х2,0 г'
$ws.Range("B57").Value = 306
$ws.Range("A58").Value = 'Фп "Щедрость природы® Фиточай при простуде" 20х2,0 г'
$ws.Range("B58").NumberFormat = "0"
$ws.Range("B58").Value = 540
$ws.Range("A59").Value = 'Фп "Щедрость природы® Фиточай кардиологический" 20х2,0 г'
$ws.Range("B59").NumberFormat = "0"
$ws.Range("B59").Value = 936
$ws.Range("A60").Value = 'Фп "Щедрость природы® Фиточай успокоительный"20х2,0 г'
$ws.Range("B60").NumberFormat = "#,##0"
$ws.Range("B60").Value = 1170
$ws.Range("A61").Value = 'Фп Сб. Грудной №4 20x2,0г'
$ws.Range("B61").Value = 458366
$ws.Range("A62").Value = 'Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г'
$ws.Range("B62").Value = 1440
$ws.Range("A63").Value = 'Фп Сб. Бруснивер 20x2,0г'
$ws.Range("B63").Value = 156894
$ws.Range("A64").Value = 'Фп Фиточай "Опалиховский" (БАД) 20х2,0 г'
$ws.Range("B64").Value = 4482
$ws.Range("A65").Value = 'Фп Фиточай "Тибетский" (БАД) 20х2,0  г'
$ws.Range("B65").Value = 8694
$ws.Range("A66").Value = 'Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г'
$ws.Range("B66").Value = 1548
$ws.Range("A67").Value = 'Фп Шалфей листья 20х1,5г'
$ws.Range("B67").Value = 131322
$ws.Range("A68").Value = 'Фп Брусника листья 20х1,5г'
$ws.Range("B68").Value = 56916
$ws.Range("A69").Value = 'Фп Череда трава 20х1,5г'
$ws.Range("B69").Value = 39635
$ws.Range("A70").Value = 'Фп Подорожник листья 20x1,5г'
$ws.Range("B70").Value = 22868
$ws.Range("A71").Value = 'Фп Пастушья сумка трава 20х1,5г'
$ws.Range("B71").Value = 4246
$ws.Range("A72").Value = 'Фп Ромашка цветки 20x1,5г'
$ws.Range("B72").Value = 1157922
$ws.Range("A73").Value = 'Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г'
$ws.Range("B73").Value = 153949
$ws.Range("A74").Value = 'Фп Сенна листья 20x1,5г'
$ws.Range("B74").Value = 64196
$ws.Range("A75").Value = 'Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г'
$ws.Range("B75").Value = 75447
$ws.Range("A76").Value = 'Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г'
$ws.Range("B76").Value = 21076
$ws.Range("A77").Value = 'Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г'
$ws.Range("B77").Value = 80761
$ws.Range("A78").Value = 'Фп Мелисса лекарственная трава 20x1,5г'
$ws.Range("B78").Value = 34974
$ws.Range("A79").Value = 'Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г'
$ws.Range("B79").Value = 82779
$ws.Range("A80").Value = 'Фп Шиповник плоды 20х2,0г'
$ws.Range("B80").Value = 47466
$ws.Range("A81").Value = 'Фп Чистотел трава 20х1,5г'
$ws.Range("B81").Value = 31380
$ws.Range("A82").Value = 'Фп Липа цветки 20x1,5г'
$ws.Range("B82").Value = 70611
$ws.Range("A83").Value = 'Фп Толокнянка листья 20x1,5г'
$ws.Range("B83").Value = 39760
$ws.Range("A84").Value = 'Фп Чабрец трава 20x1,5 г'
$ws.Range("B84").Value = 66240
$ws.Range("A85").Value = 'Фп Крапива листья 20x1,5г'
$ws.Range("B85").Value = 65063
$ws.Range("A86").Value = 'Фп Зверобой трава 20x1,5г'
$ws.Range("B86").Value = 51665
$ws.Range("A87").Value = 'Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г'
$ws.Range("B87").Value = 55459
$ws.Range("A88").Value = 'Фп Душица трава 20x1,5г'
$ws.Range("B88").Value = 29070
$ws.Range("A89").Value = 'Фп Хвощ полевой трава 20х1,5г'
$ws.Range("B89").Value = 29086
$ws.Range("A90").Value = 'Фп Пустырник трава 20x1,5г'
$ws.Range("B90").Value = 44544
$ws.Range("A91").Value = 'Фп Мята перечная листья 20x1,5г'
$ws.Range("B91").Value = 72779
$ws.Range("A92").Value = 'Фп Сб. Желудочный №3 20x2,0г'
$ws.Range("B92").Value = 24858
$ws.Range("A93").Value = 'Фп Сб. Арфазетин-Э 20x2,0г'
$ws.Range("B93").Value = 54917
$ws.Range("A94").Value = 'Фп Сб. Элекасол 20x2,0г'
$ws.Range("B94").Value = 47195
$ws.Range("A95").Value = 'Фп Фиалка трехцветная трава 20x1,5г'
$ws.Range("B95").Value = 4698
$ws.Range("A96").Value = 'Фп Береза листья 20x1,5г'
$ws.Range("B96").Value = 5746
$ws.Range("A97").Value = 'Фп Золототысячник трава 20х1,5г'
$ws.Range("B97").Value = 5367
$ws.Range("A98").Value = 'Фп Пижма цветки 20х1,5г'
$ws.Range("B98").Value = 10650
$ws.Range("A99").Value = 'Фп Боярышник плоды 20х3,0г'
$ws.Range("B99").Value = 24578
$ws.Range("A100").Value = 'Фп Аир корневища 20x1,5г'
$ws.Range("B100").Value = 6143
$ws.Range("A101").Value = 'Фп Ольха соплодия 20х1,5г'
$ws.Range("B101").Value = 5324
$ws.Range("A102").Value = 'Фп Бадан корневища 20x1,5г'
$ws.Range("B102").Value = 1249
$ws.Range("A103").Value = 'Фп Дуб кора 20х1,5г'
$ws.Range("B103").Value = 7245
$ws.Range("A104").Value = 'Фп Крушина кора 20x1,5г'
$ws.Range("B104").Value = 10226
$ws.Range("A105").Value = 'Фп Ноготки цветки 20x1,5г'
$ws.Range("B105").Value = 44123
$ws.Range("A106").Value = 'Фп Тысячелистник трава 20x1,5г'
$ws.Range("B106").Value = 26222
$ws.Range("A107").Value = 'Фп Кровохлебка корневища и корни 20x1,5г'
$ws.Range("B107").Value = 9484
$ws.Range("A108").Value = 'Фп Валериана корневища с корнями 20x1,5г'
$ws.Range("B108").Value = 28646
$ws.Range("A109").Value = 'Фп Лапчатка корневища 20x2,5г'
$ws.Range("B109").Value = 6043
$ws.Range("A110").Value = 'Фп Почечный чай листья 20x1,5г'
$ws.Range("B110").Value = 136881
$ws.Range("A111").Value = 'Фп Девясил корневища и корни 20х1,5г'
$ws.Range("B111").Value = 30255
